# Slide 1, TextBox 11 ("Assignment 2:" / topic line) - split "Assignment 2:" into
# "Assignment 2" + ":" and turn the "Assignment 2" portion into a hyperlink that
# opens an external file (the student's video for the assignment).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$para1 = $tr.Paragraphs(1)

# "Assignment 2" is the first 12 characters of the paragraph; grabbing that
# sub-range and giving it a hyperlink address causes PowerPoint to split the
# run into "Assignment 2" (with the link) and the trailing ":" (without it).
$linkChars = $para1.Characters(1, 12)
$linkChars.ActionSettings(1).Hyperlink.Address = "Assignment 2.mp4"
